$p = $ppt.ActivePresentation
try {
    $tm = $p.TitleMaster
    Write-Output ("accent1: " + $tm.ColorScheme.Colors(5).RGB)
} catch {
    Write-Output ("ERR: " + $_)
}
